$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocabulary rows: "look / appearance" (외모) topic, rows 285-294
$newWords = @(
  ,@("외모", "внешность", "look", "внешность", "noun", 2)
  ,@("예쁘다", "красивый", "look", "внешность", "adjective", 2)
  ,@("귀엽다", "милый", "look", "внешность", "adjective", 2)
  ,@("날씬하다", "стройный", "look", "внешность", "adjective", 2)
  ,@("머있다", "красивый", "look", "внешность", "adjective", 2)
  ,@("키가 크다 / 작다", "рост высокий / низкий", "look", "внешность", "adjective", 2)
  ,@("뚱뚱하다", "толстый", "look", "внешность", "adjective", 2)
  ,@("잘생기다", "красивый (хорошо появился)", "look", "внешность", "adjective", 2)
  ,@("머리가 길다 / 짧다", "волосы длинные / короткие", "look", "внешность", "adjective", 2)
  ,@("마르다", "тощий", "look", "внешность", "adjective", 2)
)

$startRow = 285
for ($i = 0; $i -lt $newWords.Count; $i++) {
  $r = $startRow + $i
  $entry = $newWords[$i]
  $ws.Cells.Item($r, 1).Value  = $entry[0]   # A: original (Korean)
  $ws.Cells.Item($r, 2).Value  = $entry[1]   # B: translate (Russian)
  $ws.Cells.Item($r, 4).Value  = $entry[2]   # D: main_thema_en
  $ws.Cells.Item($r, 5).Value  = $entry[3]   # E: main_thema_ru
  $ws.Cells.Item($r, 10).Value = $entry[4]   # J: part_of_speech
  $ws.Cells.Item($r, 11).Value = $entry[5]   # K: stage
}

# Update the sheet view: freeze header row, zoom to 55%, select J286
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 55
$ws.Range("J286").Select()
